$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume snapshot values (and a few re-ordered rows)
# Each entry: row number -> hashtable of column letter -> new text value
$updates = @{
    2 = @{ "D" = '27.103.27'; "E" = '  +0.06%  ' }
    3 = @{ "D" = '1.831.14'; "E" = '  +0.36%  ' }
    4 = @{ "D" = '1.009'; "E" = '  +0.19%  ' }
    5 = @{ "D" = '312.36'; "E" = '  +0.04%  ' }
    6 = @{ "D" = '1.008'; "E" = '  +0.18%  ' }
    7 = @{ "D" = '0.4634'; "E" = '  -1.25%  ' }
    8 = @{ "D" = '0.3714'; "E" = '  +1.66%  ' }
    9 = @{ "D" = '0.07362'; "E" = '  -0.25%  ' }
    10 = @{ "D" = '0.8738'; "E" = '  -0.41%  ' }
    11 = @{ "B" = 'TRON'; "C" = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; "D" = '0.07933'; "E" = '  +4.11%  ' }
    12 = @{ "D" = '19.92'; "E" = '  -1.60%  ' }
    13 = @{ "B" = 'WrappedEther'; "C" = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; "D" = '1.778.59'; "E" = '  -6.28%  ' }
    14 = @{ "B" = 'Chainlink'; "C" = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; "D" = '6.603'; "E" = '  +1.24%  ' }
    15 = @{ "B" = 'Polkadot'; "C" = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; "D" = '5.349'; "E" = '  -0.31%  ' }
    16 = @{ "D" = '91.98'; "E" = '  -1.39%  ' }
    17 = @{ "D" = '1.009'; "E" = '  +0.42%  ' }
    18 = @{ "D" = '0.000008895'; "E" = '  +2.11%  ' }
    19 = @{ "D" = '1.009'; "E" = '  +0.29%  ' }
    20 = @{ "B" = 'WrappedBTC'; "C" = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; "D" = '27.422.81'; "E" = '  -0.55%  ' }
    21 = @{ "B" = 'Avalanche'; "C" = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; "D" = '14.72'; "E" = '  +0.94%  ' }
    22 = @{ "D" = '5.146'; "E" = '  -1.66%  ' }
    23 = @{ "D" = '10.63'; "E" = '  +0.18%  ' }
    24 = @{ "D" = '1.999.04'; "E" = '  -4.04%  ' }
    25 = @{ "D" = '152.63'; "E" = '  +0.90%  ' }
    26 = @{ "E" = '  -1.17%  ' }
    27 = @{ "D" = '18.56'; "E" = '  +0.62%  ' }
    28 = @{ "D" = '2.095'; "E" = '  -1.50%  ' }
    29 = @{ "D" = '5.089'; "E" = '  -1.54%  ' }
    30 = @{ "D" = '115.51'; "E" = '  -0.70%  ' }
    31 = @{ "D" = '0.08871'; "E" = '  -0.48%  ' }
    32 = @{ "D" = '2.971'; "E" = '  +1.03%  ' }
    33 = @{ "D" = '0.7337'; "E" = '  -1.31%  ' }
    34 = @{ "D" = '4.457'; "E" = '  -1.11%  ' }
    35 = @{ "D" = '1.139'; "E" = '  -1.77%  ' }
    36 = @{ "D" = '2.477'; "E" = '  -8.75%  ' }
    37 = @{ "B" = 'TrustWalletToken'; "C" = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; "D" = '1.076'; "E" = '  -1.03%  ' }
    38 = @{ "B" = 'VeChain'; "C" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; "D" = '0.01950'; "E" = '  +0.93%  ' }
    39 = @{ "B" = 'Hedera'; "C" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; "D" = '0.05242'; "E" = '  -0.91%  ' }
    40 = @{ "D" = '2.935'; "E" = '  +0.33%  ' }
    41 = @{ "D" = '7.189'; "E" = '  -1.58%  ' }
    42 = @{ "D" = '0.5204'; "E" = '  -0.98%  ' }
    43 = @{ "D" = '0.8669'; "E" = '  -13.81%  ' }
    44 = @{ "D" = '0.1634'; "E" = '  -0.49%  ' }
    45 = @{ "D" = '8.239'; "E" = '  -1.46%  ' }
    46 = @{ "D" = '0.4850'; "E" = '  -0.97%  ' }
    47 = @{ "D" = '1.008'; "E" = '  +0.20%  ' }
    48 = @{ "E" = '  -1.52%  ' }
    49 = @{ "D" = '102.60'; "E" = '  -1.69%  ' }
    50 = @{ "D" = '1.631'; "E" = '  -1.14%  ' }
    51 = @{ "D" = '0.06225'; "E" = '  -0.72%  ' }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cellRef = "$col$row"
        $newValue = $updates[$row][$col]
        $cell = $ws.Range($cellRef)
        if ($col -eq "D") {
            # Column D holds price text such as "1.009" or "27.103.27".
            # Excel would silently reinterpret plain numeric-looking text as a
            # number, so force text storage with a leading apostrophe and then
            # drop the "Text" number format it applies, to keep styling identical
            # to the rest of the sheet (no explicit cell style).
            $cell.Value = "'" + $newValue
            $cell.ClearFormats()
        } else {
            $cell.Value = $newValue
        }
    }
}
